# Revert merge: restore "小组" (team/group) wording in the evaluation
# descriptions that a prior merge had changed to "学生" (student) wording,
# and move the active selection to B14 (as in the pre-merge state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C19: "对学生提交的阶段性成果进行评分" -> "对小组提交的阶段性成果进行评分"
$ws.Range("C19").Value = "对小组提交的阶段性成果进行评分"

# C20: "对学生的项目整体评分" -> "对小组的项目整体评分"
$ws.Range("C20").Value = "对小组的项目整体评分"

# Move/restore the sheet's active cell selection to B14 (was C19).
$ws.Range("B14").Select()
